$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8-30 down to 9-31
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new data record
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(8, 3).Value = "Ñuble"
$ws.Cells.Item(8, 4).Value = 45219
$ws.Cells.Item(8, 5).Value = 16
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100107
$ws.Cells.Item(8, 8).Value = "Otros"
$ws.Cells.Item(8, 9).Value = 100107002
$ws.Cells.Item(8, 10).Value = "Chirimoya"
$ws.Cells.Item(8, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 30
$ws.Cells.Item(8, 14).Value = 20000
$ws.Cells.Item(8, 15).Value = 20000
$ws.Cells.Item(8, 16).Value = 20000
$ws.Cells.Item(8, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(8, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(8, 19).Value = 2000
$ws.Cells.Item(8, 20).Value = 10

# Ensure the date style (same as other D-column cells) is applied to D8
$ws.Cells.Item(8, 4).NumberFormat = $ws.Cells.Item(9, 4).NumberFormat
